$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("U2").Value = 1.87
$ws.Range("V2").Value = 1.87
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 1.36
$ws.Range("P3").Value = 3.2
$ws.Range("U3").Value = 1.87
$ws.Range("V3").Value = 1.87
$ws.Range("V4").Value = 1.58
$ws.Range("G17").Value = 1.5
$ws.Range("H17").Value = 4.33
$ws.Range("I17").Value = 5.75
$ws.Range("J17").Value = 2.05
$ws.Range("M17").Value = 1.05
$ws.Range("N17").Value = 11
$ws.Range("Q17").Value = 1.85
$ws.Range("R17").Value = 2
$ws.Range("Z17").Value = 10
$ws.Range("AD17").Value = 8.5
$ws.Range("AG17").Value = 351
$ws.Range("AS17").Value = 151
$ws.Range("AX17").Value = 34
$ws.Range("AZ17").Value = 126
$ws.Range("G26").Value = 2.5
$ws.Range("H26").Value = 3.3
$ws.Range("I26").Value = 2.7
$ws.Range("J26").Value = 3.1
$ws.Range("N26").Value = 12
$ws.Range("W26").Value = 10
$ws.Range("Y26").Value = 10
$ws.Range("AE26").Value = 12
$ws.Range("AJ26").Value = 10
$ws.Range("AK26").Value = 26
$ws.Range("AN26").Value = 4.75
$ws.Range("AW26").Value = 4.75
$ws.Range("AZ26").Value = 41
$ws.Range("BA26").Value = 51
$ws.Range("BB26").Value = 126
$ws.Range("N37").Value = 9.85
